# Auto-generated Excel COM-interop script
# Replaces the LR-pair table (rows 2-10) with the updated TPM-derived
# 6-row table (rows 2-7): sending clusters FAPs/MuSCs x target clusters ECs/FAPs/MuSCs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old trailing rows (previously rows 8-10, MuSCs-as-sender block)
# so the table shrinks from 9 data rows to 6.
$ws.Rows("8:10").Delete() | Out-Null

# Build the replacement data block (rows 2-7, columns A-T) in one shot.
$data = New-Object 'object[,]' 6,20
# row 2
$data[0,0] = "FAPs"
$data[0,1] = "Slit1"
$data[0,2] = "Robo2"
$data[0,3] = "ECs"
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 0.2055996666666667
$data[0,7] = 0.616799
$data[0,8] = 0.9059768423248155
$data[0,9] = 0.9059768423248156
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 2.892869333333334
$data[0,13] = 8.678608
$data[0,14] = 0.9698236995656885
$data[0,15] = 0.9698236995656884
$data[0,16] = 0.5947729706435555
$data[0,17] = 5.352956735792
$data[0,18] = 0.878637812944293
$data[0,19] = 0.878637812944293
# row 3
$data[1,0] = "FAPs"
$data[1,1] = "Slit1"
$data[1,2] = "Robo2"
$data[1,3] = "FAPs"
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 0.2055996666666667
$data[1,7] = 0.616799
$data[1,8] = 0.9059768423248155
$data[1,9] = 0.9059768423248156
$data[1,10] = 1
$data[1,11] = 0.3333333333333333
$data[1,12] = 0.020314
$data[1,13] = 0.060942
$data[1,14] = 0.006810193051573731
$data[1,15] = 0.00681019305157373
$data[1,16] = 0.004176551628666667
$data[1,17] = 0.037588964658
$data[1,18] = 0.006169877196487168
$data[1,19] = 0.006169877196487168
# row 4
$data[2,0] = "FAPs"
$data[2,1] = "Slit1"
$data[2,2] = "Robo2"
$data[2,3] = "MuSCs"
$data[2,4] = 2
$data[2,5] = 0.6666666666666666
$data[2,6] = 0.2055996666666667
$data[2,7] = 0.616799
$data[2,8] = 0.9059768423248155
$data[2,9] = 0.9059768423248156
$data[2,10] = 2
$data[2,11] = 0.6666666666666666
$data[2,12] = 0.06969833333333333
$data[2,13] = 0.209095
$data[2,14] = 0.02336610738273784
$data[2,15] = 0.02336610738273783
$data[2,16] = 0.01432995410055555
$data[2,17] = 0.128969586905
$data[2,18] = 0.02116915218403538
$data[2,19] = 0.02116915218403538
# row 5
$data[3,0] = "MuSCs"
$data[3,1] = "Slit1"
$data[3,2] = "Robo2"
$data[3,3] = "ECs"
$data[3,4] = 2
$data[3,5] = 0.6666666666666666
$data[3,6] = 0.02133733333333333
$data[3,7] = 0.064012
$data[3,8] = 0.09402315767518445
$data[3,9] = 0.09402315767518446
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 2.892869333333334
$data[3,13] = 8.678608
$data[3,14] = 0.9698236995656885
$data[3,15] = 0.9698236995656884
$data[3,16] = 0.06172611725511112
$data[3,17] = 0.555535055296
$data[3,18] = 0.09118588662139544
$data[3,19] = 0.09118588662139544
# row 6
$data[4,0] = "MuSCs"
$data[4,1] = "Slit1"
$data[4,2] = "Robo2"
$data[4,3] = "FAPs"
$data[4,4] = 2
$data[4,5] = 0.6666666666666666
$data[4,6] = 0.02133733333333333
$data[4,7] = 0.064012
$data[4,8] = 0.09402315767518445
$data[4,9] = 0.09402315767518446
$data[4,10] = 1
$data[4,11] = 0.3333333333333333
$data[4,12] = 0.020314
$data[4,13] = 0.060942
$data[4,14] = 0.006810193051573731
$data[4,15] = 0.00681019305157373
$data[4,16] = 0.0004334465893333334
$data[4,17] = 0.003901019304
$data[4,18] = 0.0006403158550865623
$data[4,19] = 0.0006403158550865623
# row 7
$data[5,0] = "MuSCs"
$data[5,1] = "Slit1"
$data[5,2] = "Robo2"
$data[5,3] = "MuSCs"
$data[5,4] = 2
$data[5,5] = 0.6666666666666666
$data[5,6] = 0.02133733333333333
$data[5,7] = 0.064012
$data[5,8] = 0.09402315767518445
$data[5,9] = 0.09402315767518446
$data[5,10] = 2
$data[5,11] = 0.6666666666666666
$data[5,12] = 0.06969833333333333
$data[5,13] = 0.209095
$data[5,14] = 0.02336610738273784
$data[5,15] = 0.02336610738273783
$data[5,16] = 0.001487176571111111
$data[5,17] = 0.01338458914
$data[5,18] = 0.002196955198702451
$data[5,19] = 0.002196955198702451

$ws.Range("A2:T7").Value = $data

